# Apply BoM update: R4 and U2 were removed from the design, so:
#  - the "R2 R4" reference group becomes just "R2" (qty per group drops 2 -> 1)
#  - the "U1 U2" reference group becomes just "U1" (qty per group drops 2 -> 1)
#  - overall component counters drop by 2 (218 -> 216, 215 -> 213) on both the
#    BoM and DNF summary sheets.

$wb = $excel.ActiveWorkbook
$bom = $wb.Worksheets.Item("BoM")
$dnf = $wb.Worksheets.Item("DNF")

# --- Resistor group row (row 16): "R2 R4" -> "R2" ---
$bom.Range("E16").Value = "R2"

# --- Power supervisor group row (row 18): "U1 U2" -> "U1" ---
$bom.Range("E18").Value = "U1"

# The Quantity Per PCB / Build Quantity / Source BoM cells on these two rows
# are stored as text (shared-string) cells even though they look numeric, so
# a plain .Value assignment would coerce them to real numbers and strip the
# "t=s" shared-string typing / reset their style. Copy from existing cells
# that already hold the desired text ("1" in I13, and
# "pedalboard-display(1)" in N13) and paste values-only so the destination
# keeps its own style/format while picking up the source's text value.
$bom.Range("I13").Copy() | Out-Null
$bom.Range("I16").PasteSpecial(-4163) | Out-Null
$bom.Range("J13").Copy() | Out-Null
$bom.Range("J16").PasteSpecial(-4163) | Out-Null
$bom.Range("N13").Copy() | Out-Null
$bom.Range("N16").PasteSpecial(-4163) | Out-Null

$bom.Range("I13").Copy() | Out-Null
$bom.Range("I18").PasteSpecial(-4163) | Out-Null
$bom.Range("J13").Copy() | Out-Null
$bom.Range("J18").PasteSpecial(-4163) | Out-Null
$bom.Range("N13").Copy() | Out-Null
$bom.Range("N18").PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = $false

# --- Summary counters on both the BoM and DNF sheets ---
$bom.Range("F3").Value = "216 (205 SMD/ 0 THT)"
$bom.Range("F4").Value = "213 (204 SMD/ 0 THT)"
$bom.Range("F6").Value = 213

$dnf.Range("F3").Value = "216 (205 SMD/ 0 THT)"
$dnf.Range("F4").Value = "213 (204 SMD/ 0 THT)"
$dnf.Range("F6").Value = 213
